# Fix design age classification; add tests to validate
#
# This script reclassifies several "Screen Print Designs" rows that were
# miscategorised under the "Quick Search > New Designs" subcategory: they
# are actually older (patriotic) designs, so column F is corrected to
# "Quick Search > Patriotic" and a new "Date" column (AA) is populated to
# record when each design was actually added (1/1/2017 for the truly old
# ones, 1/1/2024 for the more recent ones) so this kind of age
# misclassification can be caught going forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Screen Print Designs")
$ws.Activate()

# --- Column F: re-point the rows that were wrongly filed under
#     "Quick Search > New Designs" to "Quick Search > Patriotic" ---
$ws.Range("F3").Value  = "Quick Search > Patriotic"
$ws.Range("F4").Value  = "Quick Search > Patriotic"
$ws.Range("F6").Value  = "Quick Search > Patriotic"
$ws.Range("F7").Value  = "Quick Search > Patriotic"
$ws.Range("F8").Value  = "Quick Search > Patriotic"
$ws.Range("F9").Value  = "Quick Search > Patriotic"
$ws.Range("F10").Value = "Quick Search > Patriotic"
$ws.Range("F11").Value = "Quick Search > Patriotic"
$ws.Range("F12").Value = "Quick Search > Patriotic"

# --- Column AA: new "Date" values used to validate the age of each design ---
$ws.Range("AA2").Value  = "1/1/2017"
$ws.Range("AA3").Value  = "1/1/2024"
$ws.Range("AA4").Value  = "1/1/2024"
$ws.Range("AA5").Value  = "1/1/2024"
$ws.Range("AA6").Value  = "1/1/2017"
$ws.Range("AA7").Value  = "1/1/2024"
$ws.Range("AA8").Value  = "1/1/2017"
$ws.Range("AA9").Value  = "1/1/2024"
$ws.Range("AA10").Value = "1/1/2024"
$ws.Range("AA11").Value = "1/1/2024"
$ws.Range("AA12").Value = "1/1/2017"
$ws.Range("AA13").Value = "1/1/2024"

# Keep the new date values stored as plain text (column AA already carries
# a text format, but make sure explicitly in case formats are recalculated).
$ws.Range("AA2:AA13").NumberFormat = "@"

# The "Embroidery Designs" sheet's selection was also moved while testing
# the same fix there.
$embroidery = $wb.Worksheets.Item("Embroidery Designs")
$embroidery.Activate()
$embroidery.Range("AA4").Select()

# "Screen Print Designs" (previously "Colors") ends up as the active /
# selected tab, with the reviewed rows selected.
$ws.Activate()
$ws.Range("F6:F12").Select()
